# The commit reorders the weekly Papaya price records (rows 2-10, excluding
# row 8 and row 11 which are unchanged) by shuffling the data columns D:T
# between rows, per the row-content permutation observed in the diff:
#   row 2 <- old row 4
#   row 3 <- old row 5
#   row 4 <- old row 2
#   row 5 <- old row 9
#   row 6 <- old row 10
#   row 7 <- old row 6
#   row 9 <- old row 7
#   row 10 <- old row 3
# (row 8 and row 11 keep their original content)
#
# Columns A-C are identical on every data row, so only D:T need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Snapshot every source row's D:T values up front (Value2 gives raw
#    numbers/strings) before any cell gets overwritten.
$row2 = $ws.Range("D2:T2").Value2
$row3 = $ws.Range("D3:T3").Value2
$row4 = $ws.Range("D4:T4").Value2
$row5 = $ws.Range("D5:T5").Value2
$row6 = $ws.Range("D6:T6").Value2
$row7 = $ws.Range("D7:T7").Value2
$row9 = $ws.Range("D9:T9").Value2
$row10 = $ws.Range("D10:T10").Value2

# 2) Write each snapshot back out to its new row.
$ws.Range("D2:T2").Value = $row4
$ws.Range("D3:T3").Value = $row5
$ws.Range("D4:T4").Value = $row2
$ws.Range("D5:T5").Value = $row9
$ws.Range("D6:T6").Value = $row10
$ws.Range("D7:T7").Value = $row6
$ws.Range("D9:T9").Value = $row7
$ws.Range("D10:T10").Value = $row3
